# Update the "取得日時" (acquisition datetime) column for the data rows
# on the "ランサーズ" sheet from "2025-10-26 01:22:50" to "2025-10-26 01:51:02".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-26 01:22:50"
$newTimestamp = "2025-10-26 01:51:02"

# Data rows are 2 through 14 (row 1 is the header row).
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
